# Update "PERIOD TO EXPIRE" (col H) and "LAST UPDATE" (col I) for rows 3-19
# to reflect new progress as of 04-Nov-2025.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

for ($row = 3; $row -le 19; $row++) {
    $hCell = $ws.Cells.Item($row, 8)   # column H: PERIOD TO EXPIRE
    $hCell.Value = $hCell.Value() - 1

    $iCell = $ws.Cells.Item($row, 9)   # column I: LAST UPDATE
    $iCell.NumberFormat = "@"          # keep as text, not auto-parsed as a date
    $iCell.Value = "04-Nov-2025"
}
